$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JasmineJourney")

$ws.Range("B4").Value = "watch me cum papi... this one is only for you 🔥"
$ws.Range("B6").Value = "I'm cumming right now... don't look away for a single second"
$ws.Range("B7").Value = "I'm so close I can feel it in every part of my body papi... wait for me, I'm right there 🔥"
$ws.Range("B8").Value = "FUCK I can't stop 🥵"
$ws.Range("B9").Value = "I've never let anyone see me like this... but you're about to 🔥"
$ws.Range("B11").Value = "watch what you're about to make me do... I can't hold it back anymore"
$ws.Range("B12").Value = "my pussy is so wet it's running down my thighs and I keep going harder and harder 🔥"
$ws.Range("B13").Value = "I'm grinding on my fingers right now imagining it's your cock and I'm losing my mind papi"
$ws.Range("B14").Value = "fuck fuck 🥵"
$ws.Range("B15").Value = "look at this papi... you did this to me and I want you to see every second 🔥"
$ws.Range("B17").Value = "what would you do to me right now if you had me? don't hold back"
$ws.Range("B18").Value = "I keep imagining you here pinning me down and it's making everything ten times more intense 🔥"
$ws.Range("B19").Value = "my fingers are already where they shouldn't be and I'm soaking wet because of you papi"
$ws.Range("B20").Value = "god... okay I wasn't expecting to feel this way 🥵"
$ws.Range("B21").Value = "this is what you're doing to me and I'm not sorry about it 🔥"
$ws.Range("B23").Value = "fuck it... I'm taking everything off and you better be ready for what's next papi 🔥"
$ws.Range("B24").Value = "something about the way you reacted just made my whole body light up... I'm getting so wet already"
$ws.Range("B25").Value = "I knew you'd like that... now I'm really starting to feel something 🥵"
